$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated for
# every data row (rows 2-66) from 45172 (2023-09-03) to 45175 (2023-09-06).
for ($row = 2; $row -le 66; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value2 = 45175
    }
}
